# Update the marksheet with corrected total marks.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: total correct-answer value changed from 3 to 5
$ws.Range("B11").Value = 5

# "Total" row: total marks changed from 51 to 85
$ws.Range("B12").Value = 85

# "Total" row, Max column: correct/total marks text changed from "45/84" to "85/140"
$ws.Range("E12").Value = "85/140"
